# "Added per line item and not line item FF and Billing reports"
#
# The per-line-item columns (Billing Cycle, Item Id, Item Name, Item Type,
# Item Unit Of measure, Item MPN, Item Period, Quantity) are being split out
# into their own report, so remove them here (old columns F:M) and leave the
# remaining "not line item" billing/customer/product/subscription columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 8 "line item" columns (old F:M): Billing Cycle, Item Id,
# Item Name, Item Type, Item Unit Of measure, Item MPN, Item Period, Quantity
$ws.Range("F1:M1").EntireColumn.Delete()

# Delta/Uom (now D:E) no longer sit next to an outlined item block - ungroup them
$ws.Range("D1:E1").EntireColumn.Ungroup()

# Refresh the AutoFilter so its range covers the new A1:Z1 extent
$ws.AutoFilterMode = $false
$ws.Range("A1:Z1").AutoFilter()

# Refresh the hidden _FilterDatabase defined name to match
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$Z`$1"
    }
}

# Reset the view: scroll back to column A and select E2 (was scrolled to S1 /
# W2 selected when the sheet had more columns)
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E2").Select()
